$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-09 06:48:41'
$ws.Range('O2').Value = '-3.9 °C'
$ws.Range('E3').Value = '2026-02-09 06:48:44'
$ws.Range('O3').Value = '-6.0 °C'
$ws.Range('E4').Value = '2026-02-09 06:48:46'
$ws.Range('O4').Value = '4.1 °C'
$ws.Range('E5').Value = '2026-02-09 06:48:49'
$ws.Range('M5').Value = '-3.8 °C 6:17 TU'
$ws.Range('O5').Value = '-5.3 °C'
$ws.Range('E6').Value = '2026-02-09 06:48:51'
$ws.Range('O6').Value = '6.2 °C'
$ws.Range('E7').Value = '2026-02-09 06:48:53'
$ws.Range('H7').Value = "'69%"
$ws.Range('L7').Value = '14.0 km/h - 320º 6:12 TU'
$ws.Range('E8').Value = '2026-02-09 06:48:56'
$ws.Range('L8').Value = '23.0 km/h - 309º 6:24 TU'
$ws.Range('N8').Value = '6.5 °C 6:29 TU'
$ws.Range('O8').Value = '7.1 °C'
$ws.Range('E9').Value = '2026-02-09 06:48:58'
$ws.Range('H9').Value = "'90%"
$ws.Range('N9').Value = '2.3 °C 6:28 TU'
$ws.Range('O9').Value = '6.1 °C'
$ws.Range('E10').Value = '2026-02-09 06:49:01'
$ws.Range('O10').Value = '5.0 °C'
$ws.Range('E11').Value = '2026-02-09 06:49:03'
$ws.Range('E12').Value = '2026-02-09 06:49:06'
$ws.Range('H12').Value = "'93%"
$ws.Range('N12').Value = '2.9 °C 6:21 TU'
$ws.Range('O12').Value = '6.7 °C'
$ws.Range('E13').Value = '2026-02-09 06:49:08'
$ws.Range('E14').Value = '2026-02-09 06:49:11'
$ws.Range('H14').Value = "'95%"
$ws.Range('L14').Value = '18.4 km/h - 310º 6:11 TU'
$ws.Range('E15').Value = '2026-02-09 06:49:13'
$ws.Range('H15').Value = "'89%"
$ws.Range('O15').Value = '5.0 °C'
$ws.Range('E16').Value = '2026-02-09 06:49:15'
$ws.Range('H16').Value = "'64%"
$ws.Range('O16').Value = '-5.4 °C'
$ws.Range('E17').Value = '2026-02-09 06:49:18'
$ws.Range('L17').Value = '35.3 km/h - 257º 6:25 TU'
$ws.Range('N17').Value = '-1.2 °C 6:29 TU'
$ws.Range('O17').Value = '-0.2 °C'
$ws.Range('E18').Value = '2026-02-09 06:49:21'
$ws.Range('N18').Value = '3.1 °C 6:28 TU'
$ws.Range('O18').Value = '5.7 °C'
$ws.Range('E19').Value = '2026-02-09 06:49:23'
$ws.Range('N19').Value = '2.3 °C 6:26 TU'
$ws.Range('O19').Value = '3.1 °C'
$ws.Range('E20').Value = '2026-02-09 06:49:25'
$ws.Range('M20').Value = '-4.9 °C 6:19 TU'
$ws.Range('O20').Value = '-6.2 °C'
$ws.Range('E21').Value = '2026-02-09 06:49:28'
$ws.Range('J21').Value = '1010.0 hPa'
$ws.Range('O21').Value = '0.4 °C'
$ws.Range('E22').Value = '2026-02-09 06:49:30'
$ws.Range('E23').Value = '2026-02-09 06:49:33'
$ws.Range('E24').Value = '2026-02-09 06:49:36'
$ws.Range('H24').Value = "'90%"
$ws.Range('O24').Value = '4.4 °C'
$ws.Range('E25').Value = '2026-02-09 06:49:38'
$ws.Range('L25').Value = '22.7 km/h - 246º 6:23 TU'
$ws.Range('O25').Value = '-4.2 °C'
$ws.Range('E26').Value = '2026-02-09 06:49:41'
$ws.Range('J26').Value = '1008.6 hPa'
$ws.Range('E27').Value = '2026-02-09 06:49:44'
$ws.Range('E28').Value = '2026-02-09 06:49:46'
$ws.Range('E29').Value = '2026-02-09 06:49:48'
$ws.Range('N29').Value = '2.6 °C 6:08 TU'
$ws.Range('O29').Value = '4.9 °C'
$ws.Range('E30').Value = '2026-02-09 06:49:51'
$ws.Range('J30').Value = '1007.7 hPa'
$ws.Range('N30').Value = '4.4 °C 6:14 TU'
$ws.Range('O30').Value = '6.2 °C'
$ws.Range('E31').Value = '2026-02-09 06:49:53'
$ws.Range('J31').Value = '1006.7 hPa'
$ws.Range('O31').Value = '8.8 °C'
$ws.Range('E32').Value = '2026-02-09 06:49:56'
$ws.Range('H32').Value = "'79%"
$ws.Range('K32').Value = '-0.1 MJ/m2'
$ws.Range('E33').Value = '2026-02-09 06:49:59'
$ws.Range('H33').Value = "'94%"
$ws.Range('J33').Value = '1009.8 hPa'
$ws.Range('E34').Value = '2026-02-09 06:50:02'
$ws.Range('E35').Value = '2026-02-09 06:50:04'
$ws.Range('H35').Value = "'67%"
$ws.Range('N35').Value = '3.0 °C 6:28 TU'
$ws.Range('E36').Value = '2026-02-09 06:50:07'
$ws.Range('H36').Value = "'86%"
$ws.Range('O36').Value = '7.8 °C'
$ws.Range('E37').Value = '2026-02-09 06:50:10'
$ws.Range('O37').Value = '2.9 °C'
$ws.Range('E38').Value = '2026-02-09 06:50:12'
$ws.Range('N38').Value = '3.2 °C 6:23 TU'
$ws.Range('O38').Value = '5.6 °C'
$ws.Range('E39').Value = '2026-02-09 06:50:15'
$ws.Range('E40').Value = '2026-02-09 06:50:18'
$ws.Range('O40').Value = '-0.2 °C'
$ws.Range('E41').Value = '2026-02-09 06:50:20'
$ws.Range('H41').Value = "'55%"
$ws.Range('E42').Value = '2026-02-09 06:50:23'
$ws.Range('H42').Value = "'98%"
$ws.Range('N42').Value = '3.6 °C 6:26 TU'
$ws.Range('O42').Value = '5.9 °C'
$ws.Range('E43').Value = '2026-02-09 06:50:25'
$ws.Range('N43').Value = '5.5 °C 6:24 TU'
$ws.Range('E44').Value = '2026-02-09 06:50:28'
$ws.Range('O44').Value = '-7.0 °C'
$ws.Range('E45').Value = '2026-02-09 06:50:31'
$ws.Range('J45').Value = '1009.4 hPa'
$ws.Range('L45').Value = '10.1 km/h - 170º 6:19 TU'
$ws.Range('M45').Value = '1.8 °C 6:22 TU'
$ws.Range('O45').Value = '0.2 °C'
$ws.Range('E46').Value = '2026-02-09 06:50:34'
$ws.Range('J46').Value = '1009.5 hPa'
$ws.Range('K46').Value = '-0.1 MJ/m2'
